$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Interpolation_instructions")
$ws2 = $wb.Worksheets.Item("Trend_instructions")

# --- Trend_instructions sheet ---
# Row 2: end_year 1936 -> 1934
$ws2.Range("E2").Value = 1934

# Row 3: CEDS_fuel coal_coke -> hard_coal, CEDS_sector autoproducer -> public, start_year 1931 -> 1932
$ws2.Range("B3").Value = "hard_coal"
$ws2.Range("C3").Value = "1A1a_Electricity-public"
$ws2.Range("D3").Value = 1932

# Row 4: CEDS_fuel coal_coke -> brown_coal, CEDS_sector Heat-production -> Electricity-public
$ws2.Range("B4").Value = "brown_coal"
$ws2.Range("C4").Value = "1A1a_Electricity-public"

# Row 5: CEDS_sector Other-transformation -> Electricity-autoproducer
$ws2.Range("C5").Value = "1A1a_Electricity-autoproducer"

# --- Interpolation_instructions sheet ---
# method for row 2 changes from "match_to_trend" to "linear"
$ws1.Range("H2").Value = "linear"
$ws1.Range("H2").Select()

# Widen column C to fit its new longer text values
$ws2.Columns.Item(3).ColumnWidth = 25.166666666666668

# Make Trend_instructions the active sheet/tab with H12 selected
$ws2.Activate()
$ws2.Range("H12").Select()
